# aggiornamento fino a 20/09/2021
# Appends daily rows 375-385 (dates 2021-09-10 .. 2021-09-20) to Sheet1,
# continuing the existing A:D table (date serial, nuovi pos., somma mobile
# 7gg., somma mobile 7gg. per 100mila abitanti).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting (style s="2": bold, centered, bordered,
# YYYY-MM-DD HH:MM:SS number format) from the last existing row (A374)
# down onto the new date cells before writing values.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)

$newRows = @(
    @(375, 44449, 0, 0, 0),
    @(376, 44450, 0, 0, 0),
    @(377, 44451, 1, 1, 46.70714619336758),
    @(378, 44452, 0, 1, 46.70714619336758),
    @(379, 44453, 0, 1, 46.70714619336758),
    @(380, 44454, 0, 1, 46.70714619336758),
    @(381, 44455, 0, 1, 46.70714619336758),
    @(382, 44456, 0, 1, 46.70714619336758),
    @(383, 44457, 0, 1, 46.70714619336758),
    @(384, 44458, 0, 0, 0),
    @(385, 44459, 0, 0, 0)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

Write-Output "Appended rows 375-385 (through 2021-09-20)"
